$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87, shifting existing rows 87:166 down to 88:167
$ws.Rows("87:87").Insert()

# Populate the newly inserted row 87 with its data
$ws.Range("A87").Value = 9
$ws.Range("B87").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C87").Value = "Metropolitana"
$ws.Range("D87").Value = 44484
$ws.Range("E87").Value = 13
$ws.Range("F87").Value = 300000001
$ws.Range("G87").Value = "Rabanito"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 7900
$ws.Range("K87").Value = 3000
$ws.Range("L87").Value = 4000
$ws.Range("M87").Value = 3494
$ws.Range("N87").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O87").Value = "Provincia de Chacabuco"
$ws.Range("P87").Value = 35
$ws.Range("Q87").Value = 100
$ws.Range("R87").Value = "Hortaliza"
